# Adding title name at top of nav. Prepping for css grid
#
# 1) Breakpoint Calculations: move the selection to G8 (was B11).
# 2) Add a new "Grid_Mobile" worksheet after "Skills" with the column /
#    margin / gutter grid labels + dimensions, and make it the active tab.

$wb = $excel.ActiveWorkbook

# --- 1. Update selection on "Breakpoint Calculations" -----------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("G8").Select()

# --- 2. Add the new "Grid_Mobile" sheet after the last existing sheet -
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Grid_Mobile"

# Write the header row / dimension row. The order below matches the
# order the new labels were first typed in (so shared-string indices
# line up with the authored workbook).
$newSheet.Range("C1").Value = "Column 1"
$newSheet.Range("E1").Value = "Column 2"
$newSheet.Range("G1").Value = "Column 3"
$newSheet.Range("I1").Value = "Column 4"
$newSheet.Range("B1").Value = "Margin"
$newSheet.Range("D1").Value = "Gutter"
$newSheet.Range("C2").Value = "auto"
$newSheet.Range("A2").Value = "Dimensions"
$newSheet.Range("B2").Value = "16px"

$newSheet.Range("F1").Value = "Gutter"
$newSheet.Range("H1").Value = "Gutter"
$newSheet.Range("J1").Value = "Margin"
$newSheet.Range("D2").Value = "16px"
$newSheet.Range("E2").Value = "auto"
$newSheet.Range("F2").Value = "16px"
$newSheet.Range("G2").Value = "auto"
$newSheet.Range("H2").Value = "16px"
$newSheet.Range("I2").Value = "auto"
$newSheet.Range("J2").Value = "16px"

# Size column A to fit the "Dimensions" label.
$newSheet.Columns.Item(1).AutoFit() | Out-Null

# Leave the new sheet selected on G20, matching the authored file.
$newSheet.Range("G20").Select()
